$d = $word.ActiveDocument

$d.Content.Find.Execute("39-28=11", $true, $false, $false, $false, $false, $true, 1, $false, "88-27=61", 1) | Out-Null
$d.Content.Find.Execute("34-26=8", $true, $false, $false, $false, $false, $true, 1, $false, "59-41=18", 1) | Out-Null
$d.Content.Find.Execute("16+8=24", $true, $false, $false, $false, $false, $true, 1, $false, "53-37=16", 1) | Out-Null
$d.Content.Find.Execute("68+17=85", $true, $false, $false, $false, $false, $true, 1, $false, "87-2=85", 1) | Out-Null
$d.Content.Find.Execute("90-32=58", $true, $false, $false, $false, $false, $true, 1, $false, "71-18=53", 1) | Out-Null
$d.Content.Find.Execute("34+15=49", $true, $false, $false, $false, $false, $true, 1, $false, "35+14=49", 1) | Out-Null
$d.Content.Find.Execute("51+34=85", $true, $false, $false, $false, $false, $true, 1, $false, "15+41=56", 1) | Out-Null
$d.Content.Find.Execute("0+38=38", $true, $false, $false, $false, $false, $true, 1, $false, "28-4=24", 1) | Out-Null
$d.Content.Find.Execute("54-52=2", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=52", 1) | Out-Null
$d.Content.Find.Execute("44-22=22", $true, $false, $false, $false, $false, $true, 1, $false, "59-44=15", 1) | Out-Null
$d.Content.Find.Execute("56+20=76", $true, $false, $false, $false, $false, $true, 1, $false, "1+36=37", 1) | Out-Null
$d.Content.Find.Execute("47-7=40", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=37", 1) | Out-Null
$d.Content.Find.Execute("83-77=6", $true, $false, $false, $false, $false, $true, 1, $false, "12+1=13", 1) | Out-Null
$d.Content.Find.Execute("90-14=76", $true, $false, $false, $false, $false, $true, 1, $false, "51-38=13", 1) | Out-Null
$d.Content.Find.Execute("16-5=11", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=71", 1) | Out-Null
$d.Content.Find.Execute("16-13=3", $true, $false, $false, $false, $false, $true, 1, $false, "74+3=77", 1) | Out-Null
$d.Content.Find.Execute("41+52=93", $true, $false, $false, $false, $false, $true, 1, $false, "38+54=92", 1) | Out-Null
$d.Content.Find.Execute("52+31=83", $true, $false, $false, $false, $false, $true, 1, $false, "88-62=26", 1) | Out-Null
$d.Content.Find.Execute("59-42=17", $true, $false, $false, $false, $false, $true, 1, $false, "52-36=16", 1) | Out-Null
$d.Content.Find.Execute("55-1=54", $true, $false, $false, $false, $false, $true, 1, $false, "16+31=47", 1) | Out-Null
$d.Content.Find.Execute("96-84=12", $true, $false, $false, $false, $false, $true, 1, $false, "88-46=42", 1) | Out-Null
$d.Content.Find.Execute("78-61=17", $true, $false, $false, $false, $false, $true, 1, $false, "65-2=63", 1) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "22+11=33", 1) | Out-Null
$d.Content.Find.Execute("91-25=66", $true, $false, $false, $false, $false, $true, 1, $false, "58-26=32", 1) | Out-Null
$d.Content.Find.Execute("94-37=57", $true, $false, $false, $false, $false, $true, 1, $false, "99-21=78", 1) | Out-Null
$d.Content.Find.Execute("62-12=50", $true, $false, $false, $false, $false, $true, 1, $false, "73-72=1", 1) | Out-Null
$d.Content.Find.Execute("83-34=49", $true, $false, $false, $false, $false, $true, 1, $false, "83+8=91", 1) | Out-Null
$d.Content.Find.Execute("30+15=45", $true, $false, $false, $false, $false, $true, 1, $false, "42-7=35", 1) | Out-Null
$d.Content.Find.Execute("67-26=41", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=79", 1) | Out-Null
$d.Content.Find.Execute("90-87=3", $true, $false, $false, $false, $false, $true, 1, $false, "61-13=48", 1) | Out-Null
$d.Content.Find.Execute("70-60=10", $true, $false, $false, $false, $false, $true, 1, $false, "38+51=89", 1) | Out-Null
$d.Content.Find.Execute("90-5=85", $true, $false, $false, $false, $false, $true, 1, $false, "40+23=63", 1) | Out-Null
$d.Content.Find.Execute("83-5=78", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=34", 1) | Out-Null
$d.Content.Find.Execute("15+31=46", $true, $false, $false, $false, $false, $true, 1, $false, "77-44=33", 1) | Out-Null
$d.Content.Find.Execute("40-39=1", $true, $false, $false, $false, $false, $true, 1, $false, "22-8=14", 1) | Out-Null
$d.Content.Find.Execute("67-42=25", $true, $false, $false, $false, $false, $true, 1, $false, "23-10=13", 1) | Out-Null
$d.Content.Find.Execute("90-4=86", $true, $false, $false, $false, $false, $true, 1, $false, "49+16=65", 1) | Out-Null
$d.Content.Find.Execute("85+3=88", $true, $false, $false, $false, $false, $true, 1, $false, "50+9=59", 1) | Out-Null
$d.Content.Find.Execute("89+6=95", $true, $false, $false, $false, $false, $true, 1, $false, "16+48=64", 1) | Out-Null
$d.Content.Find.Execute("47+17=64", $true, $false, $false, $false, $false, $true, 1, $false, "36+58=94", 1) | Out-Null
$d.Content.Find.Execute("58+5=63", $true, $false, $false, $false, $false, $true, 1, $false, "93-34=59", 1) | Out-Null
$d.Content.Find.Execute("27+17=44", $true, $false, $false, $false, $false, $true, 1, $false, "87-71=16", 1) | Out-Null
$d.Content.Find.Execute("19+60=79", $true, $false, $false, $false, $false, $true, 1, $false, "82-19=63", 1) | Out-Null
$d.Content.Find.Execute("83-69=14", $true, $false, $false, $false, $false, $true, 1, $false, "52+38=90", 1) | Out-Null
$d.Content.Find.Execute("83-13=70", $true, $false, $false, $false, $false, $true, 1, $false, "9+65=74", 1) | Out-Null
$d.Content.Find.Execute("53-13=40", $true, $false, $false, $false, $false, $true, 1, $false, "39+54=93", 1) | Out-Null
$d.Content.Find.Execute("39+45=84", $true, $false, $false, $false, $false, $true, 1, $false, "81+14=95", 1) | Out-Null
$d.Content.Find.Execute("0+18=18", $true, $false, $false, $false, $false, $true, 1, $false, "98-9=89", 1) | Out-Null
$d.Content.Find.Execute("81-57=24", $true, $false, $false, $false, $false, $true, 1, $false, "16-1=15", 1) | Out-Null
$d.Content.Find.Execute("76-12=64", $true, $false, $false, $false, $false, $true, 1, $false, "47+16=63", 1) | Out-Null
$d.Content.Find.Execute("4+15=19", $true, $false, $false, $false, $false, $true, 1, $false, "80-61=19", 1) | Out-Null
$d.Content.Find.Execute("21+61=82", $true, $false, $false, $false, $false, $true, 1, $false, "73-36=37", 1) | Out-Null
$d.Content.Find.Execute("23+26=49", $true, $false, $false, $false, $false, $true, 1, $false, "31-18=13", 1) | Out-Null
$d.Content.Find.Execute("53-47=6", $true, $false, $false, $false, $false, $true, 1, $false, "7+1=8", 1) | Out-Null
$d.Content.Find.Execute("64-14=50", $true, $false, $false, $false, $false, $true, 1, $false, "82-20=62", 1) | Out-Null
$d.Content.Find.Execute("72+13=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+36=50", 1) | Out-Null
$d.Content.Find.Execute("24-18=6", $true, $false, $false, $false, $false, $true, 1, $false, "15+12=27", 1) | Out-Null
$d.Content.Find.Execute("55-26=29", $true, $false, $false, $false, $false, $true, 1, $false, "5+53=58", 1) | Out-Null
$d.Content.Find.Execute("35-29=6", $true, $false, $false, $false, $false, $true, 1, $false, "85+5=90", 1) | Out-Null
$d.Content.Find.Execute("61-53=8", $true, $false, $false, $false, $false, $true, 1, $false, "54+8=62", 1) | Out-Null
$d.Content.Find.Execute("13-3=10", $true, $false, $false, $false, $false, $true, 1, $false, "13+84=97", 1) | Out-Null
$d.Content.Find.Execute("16+26=42", $true, $false, $false, $false, $false, $true, 1, $false, "99-46=53", 1) | Out-Null
$d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "99-18=81", 1) | Out-Null
$d.Content.Find.Execute("30+67=97", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=20", 1) | Out-Null
$d.Content.Find.Execute("80-76=4", $true, $false, $false, $false, $false, $true, 1, $false, "3+30=33", 1) | Out-Null
$d.Content.Find.Execute("20+23=43", $true, $false, $false, $false, $false, $true, 1, $false, "49+38=87", 1) | Out-Null
$d.Content.Find.Execute("85-49=36", $true, $false, $false, $false, $false, $true, 1, $false, "20+65=85", 1) | Out-Null
$d.Content.Find.Execute("3+28=31", $true, $false, $false, $false, $false, $true, 1, $false, "7+25=32", 1) | Out-Null
$d.Content.Find.Execute("95-92=3", $true, $false, $false, $false, $false, $true, 1, $false, "71-2=69", 1) | Out-Null
$d.Content.Find.Execute("29+0=29", $true, $false, $false, $false, $false, $true, 1, $false, "24-11=13", 1) | Out-Null
$d.Content.Find.Execute("71-49=22", $true, $false, $false, $false, $false, $true, 1, $false, "69-68=1", 1) | Out-Null
$d.Content.Find.Execute("10+59=69", $true, $false, $false, $false, $false, $true, 1, $false, "9+30=39", 1) | Out-Null
$d.Content.Find.Execute("57+25=82", $true, $false, $false, $false, $false, $true, 1, $false, "12+25=37", 1) | Out-Null
$d.Content.Find.Execute("31+11=42", $true, $false, $false, $false, $false, $true, 1, $false, "13+75=88", 1) | Out-Null
$d.Content.Find.Execute("50-12=38", $true, $false, $false, $false, $false, $true, 1, $false, "0+37=37", 1) | Out-Null
$d.Content.Find.Execute("36-34=2", $true, $false, $false, $false, $false, $true, 1, $false, "82-61=21", 1) | Out-Null
$d.Content.Find.Execute("45-16=29", $true, $false, $false, $false, $false, $true, 1, $false, "88-16=72", 1) | Out-Null
$d.Content.Find.Execute("14+42=56", $true, $false, $false, $false, $false, $true, 1, $false, "84-40=44", 1) | Out-Null
$d.Content.Find.Execute("25+6=31", $true, $false, $false, $false, $false, $true, 1, $false, "18-17=1", 1) | Out-Null
$d.Content.Find.Execute("27-17=10", $true, $false, $false, $false, $false, $true, 1, $false, "75-49=26", 1) | Out-Null
$d.Content.Find.Execute("39+16=55", $true, $false, $false, $false, $false, $true, 1, $false, "19+20=39", 1) | Out-Null
$d.Content.Find.Execute("33+35=68", $true, $false, $false, $false, $false, $true, 1, $false, "93+2=95", 1) | Out-Null
$d.Content.Find.Execute("18+35=53", $true, $false, $false, $false, $false, $true, 1, $false, "94-86=8", 1) | Out-Null
$d.Content.Find.Execute("58+29=87", $true, $false, $false, $false, $false, $true, 1, $false, "55-23=32", 1) | Out-Null
$d.Content.Find.Execute("1+50=51", $true, $false, $false, $false, $false, $true, 1, $false, "30+9=39", 1) | Out-Null
$d.Content.Find.Execute("66+13=79", $true, $false, $false, $false, $false, $true, 1, $false, "35+41=76", 1) | Out-Null
$d.Content.Find.Execute("34-6=28", $true, $false, $false, $false, $false, $true, 1, $false, "50+30=80", 1) | Out-Null
$d.Content.Find.Execute("61-50=11", $true, $false, $false, $false, $false, $true, 1, $false, "98+0=98", 1) | Out-Null
$d.Content.Find.Execute("20+8=28", $true, $false, $false, $false, $false, $true, 1, $false, "76+6=82", 1) | Out-Null
$d.Content.Find.Execute("38-35=3", $true, $false, $false, $false, $false, $true, 1, $false, "65+22=87", 1) | Out-Null
$d.Content.Find.Execute("46+37=83", $true, $false, $false, $false, $false, $true, 1, $false, "62-3=59", 1) | Out-Null
$d.Content.Find.Execute("46-25=21", $true, $false, $false, $false, $false, $true, 1, $false, "93-90=3", 1) | Out-Null
$d.Content.Find.Execute("78-22=56", $true, $false, $false, $false, $false, $true, 1, $false, "30+25=55", 1) | Out-Null
$d.Content.Find.Execute("19+62=81", $true, $false, $false, $false, $false, $true, 1, $false, "83-76=7", 1) | Out-Null
$d.Content.Find.Execute("47+30=77", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=34", 1) | Out-Null
$d.Content.Find.Execute("48+10=58", $true, $false, $false, $false, $false, $true, 1, $false, "82+17=99", 1) | Out-Null
$d.Content.Find.Execute("44+20=64", $true, $false, $false, $false, $false, $true, 1, $false, "98-46=52", 1) | Out-Null
$d.Content.Find.Execute("96-52=44", $true, $false, $false, $false, $false, $true, 1, $false, "10+10=20", 1) | Out-Null
$d.Content.Find.Execute("21+61=82", $true, $false, $false, $false, $false, $true, 1, $false, "61+28=89", 1) | Out-Null
$d.Content.Find.Execute("5+82=87", $true, $false, $false, $false, $false, $true, 1, $false, "59+20=79", 1) | Out-Null
